$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Update the "Syllabus version" timestamp in the body text.
# -----------------------------------------------------------------
$d.Content.Find.Execute(
    "Syllabus version: 05:11 PM, 12 February, 2018", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Syllabus version: 05:13 PM, 12 February, 2018", 2) | Out-Null

# -----------------------------------------------------------------
# 2) Insert a new "Exercise:" bullet (with a link to the discrete
#    characters exercise repo) right before the "Week 7, Feb 20 &
#    Feb 22" heading, as a sibling of the other Week 6 bullets
#    (Topic / Items / Video), reusing their numId (1015) at ilvl 1.
# -----------------------------------------------------------------
$week7 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Week 7, Feb 20*") {
        $week7 = $p
        break
    }
}

if ($week7 -ne $null) {
    # Create a blank paragraph right before "Week 7, Feb 20 & Feb 22"
    $week7.Range.InsertParagraphBefore()

    # Re-locate "Week 7..." (paragraph collection is now shifted) and
    # grab the new blank paragraph that now sits right before it.
    $week7 = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like "*Week 7, Feb 20*") {
            $week7 = $p
            break
        }
    }
    $newPara = $week7.Previous()

    # Replace the whole blank paragraph (including its end-of-paragraph
    # mark) with the fully-formed "Exercise:" paragraph, via OOXML, so
    # we get precise control of style / numbering / hyperlink markup.
    $insertRange = $d.Range($newPara.Range.Start, $newPara.Range.End)

    $xmlFrag = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:pStyle w:val="Compact"/><w:numPr><w:numId w:val="1015"/><w:ilvl w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Exercise:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:hyperlink r:id="rIdNEWEXERCISE"><w:r><w:t xml:space="preserve">https://github.com/bomeara/phylometh_discrete</w:t></w:r></w:hyperlink></w:p></w:body></w:document></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/_rels/document.xml.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rIdNEWEXERCISE" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink" Target="https://github.com/bomeara/phylometh_discrete" TargetMode="External"/></Relationships></pkg:xmlData></pkg:part></pkg:package>
'@
    $insertRange.InsertXML($xmlFrag)

    # The hyperlink run was imported without its "Hyperlink" character
    # style (a quirk of raw XML import), so re-apply it explicitly by
    # locating the newly created hyperlink through the Hyperlinks
    # collection and restyling its Range.
    foreach ($h in $d.Hyperlinks) {
        if ($h.Address -eq "https://github.com/bomeara/phylometh_discrete") {
            $h.Range.Style = "Hyperlink"
        }
    }
}

$d.Saved = $false
